$wb = $excel.ActiveWorkbook

# Sheet 1
$ws = $wb.Worksheets.Item(1)
$ws.Range("D7").Value = 0.03098275708614396
$ws.Range("E7").Value = 0.1345085345116891
$ws.Range("F7").Value = 0.1797990992421733
$ws.Range("G7").Value = 0.176019195220703
$ws.Range("H7").Value = 19.46300214859709
$ws.Range("L7").Value = 0.03117394853389245
$ws.Range("M7").Value = 0.1419980679078272
$ws.Range("N7").Value = 0.1973562431169176
$ws.Range("O7").Value = 0.1765614582344982
$ws.Range("P7").Value = 20.89496329383445

# Sheet 2
$ws = $wb.Worksheets.Item(2)
$ws.Range("D7").Value = 0.03273117557743174
$ws.Range("E7").Value = 0.1257064520503097
$ws.Range("F7").Value = 0.2175013606752148
$ws.Range("G7").Value = 0.1809175933330745
$ws.Range("H7").Value = 21.69918585846215
$ws.Range("L7").Value = 0.03056095591397193
$ws.Range("M7").Value = 0.1165397114344599
$ws.Range("N7").Value = 0.1790653980621287
$ws.Range("O7").Value = 0.1748169211317141
$ws.Range("P7").Value = 18.7322868933271

# Sheet 3
$ws = $wb.Worksheets.Item(3)
$ws.Range("D7").Value = 0.03059809449925932
$ws.Range("E7").Value = 0.1266480691549654
$ws.Range("F7").Value = 0.2032830012941191
$ws.Range("G7").Value = 0.1749231102492158
$ws.Range("H7").Value = 20.63755080282857
$ws.Range("L7").Value = 0.02615616805087571
$ws.Range("M7").Value = 0.1139688665360085
$ws.Range("N7").Value = 0.1830358624993586
$ws.Range("O7").Value = 0.1617286865428509
$ws.Range("P7").Value = 18.47503400918239

# Sheet 4
$ws = $wb.Worksheets.Item(4)
$ws.Range("D7").Value = 0.1844804259461474
$ws.Range("E7").Value = 0.3486379369538254
$ws.Range("F7").Value = 0.1244832377305674
$ws.Range("G7").Value = 0.4295118461068884
$ws.Range("H7").Value = 13.24314325736896
$ws.Range("L7").Value = 0.1786813513902775
$ws.Range("M7").Value = 0.3413369000438148
$ws.Range("N7").Value = 0.1223146258380136
$ws.Range("O7").Value = 0.4227071697881141
$ws.Range("P7").Value = 13.4363749652801

# Sheet 5
$ws = $wb.Worksheets.Item(5)
$ws.Range("D7").Value = 0.3161936410902922
$ws.Range("E7").Value = 0.4316364297772007
$ws.Range("F7").Value = 0.2201515474748739
$ws.Range("G7").Value = 0.5623109825446166
$ws.Range("H7").Value = 18.77027984749002
$ws.Range("L7").Value = 0.2224879338492728
$ws.Range("M7").Value = 0.3546982246546929
$ws.Range("N7").Value = 0.1518178089784719
$ws.Range("O7").Value = 0.4716862663352335
$ws.Range("P7").Value = 15.20431011501121

# Sheet 6
$ws = $wb.Worksheets.Item(6)
$ws.Range("D7").Value = 0.1867499245041501
$ws.Range("E7").Value = 0.3152026029584383
$ws.Range("F7").Value = 0.1210780969038256
$ws.Range("G7").Value = 0.4321457213766557
$ws.Range("H7").Value = 12.77755325252691
$ws.Range("L7").Value = 0.2154415627506814
$ws.Range("M7").Value = 0.3657175233533907
$ws.Range("N7").Value = 0.1393619083607081
$ws.Range("O7").Value = 0.4641568299084711
$ws.Range("P7").Value = 14.8487868294915

# Sheet 7
$ws = $wb.Worksheets.Item(7)
$ws.Range("D7").Value = 0.03666660172407907
$ws.Range("E7").Value = 0.1630022134702198
$ws.Range("F7").Value = 0.2308754647560251
$ws.Range("G7").Value = 0.1914852519753912
$ws.Range("H7").Value = 24.53505366686255
$ws.Range("L7").Value = 0.03406749380681216
$ws.Range("M7").Value = 0.146950580273663
$ws.Range("N7").Value = 0.205139490666727
$ws.Range("O7").Value = 0.1845738166880995
$ws.Range("P7").Value = 21.97531281743247

# Sheet 8
$ws = $wb.Worksheets.Item(8)
$ws.Range("D7").Value = 0.04113831825349686
$ws.Range("E7").Value = 0.1377801668787181
$ws.Range("F7").Value = 0.2742471685402062
$ws.Range("G7").Value = 0.2028258323130879
$ws.Range("H7").Value = 25.09454205865735
$ws.Range("L7").Value = 0.04358375156416795
$ws.Range("M7").Value = 0.1368005707400476
$ws.Range("N7").Value = 0.2907827019720263
$ws.Range("O7").Value = 0.2087672186052397
$ws.Range("P7").Value = 24.91057501013136

# Sheet 9
$ws = $wb.Worksheets.Item(9)
$ws.Range("D7").Value = 0.0381418034392564
$ws.Range("E7").Value = 0.1374068262547268
$ws.Range("F7").Value = 0.2285044802011466
$ws.Range("G7").Value = 0.1952992663561653
$ws.Range("H7").Value = 22.5907030955483
$ws.Range("L7").Value = 0.03579435691763709
$ws.Range("M7").Value = 0.1332000726735864
$ws.Range("N7").Value = 0.2163857600234187
$ws.Range("O7").Value = 0.1891939663880355
$ws.Range("P7").Value = 21.73922264810598

# Sheet 10
$ws = $wb.Worksheets.Item(10)
$ws.Range("D7").Value = 0.563989668305375
$ws.Range("E7").Value = 0.6050020707591912
$ws.Range("F7").Value = 0.250376755295017
$ws.Range("G7").Value = 0.7509924555582266
$ws.Range("H7").Value = 26.52580697347074
$ws.Range("L7").Value = 0.6025211658214754
$ws.Range("M7").Value = 0.6066538112178462
$ws.Range("N7").Value = 0.2556503810929618
$ws.Range("O7").Value = 0.7762223687974183
$ws.Range("P7").Value = 27.14599183009396

# Sheet 11
$ws = $wb.Worksheets.Item(11)
$ws.Range("D7").Value = 0.3669741446728907
$ws.Range("E7").Value = 0.4660231200039454
$ws.Range("F7").Value = 0.1727462728139513
$ws.Range("G7").Value = 0.6057839092224971
$ws.Range("H7").Value = 18.02691542229411
$ws.Range("L7").Value = 0.4835060654124726
$ws.Range("M7").Value = 0.4897680301519226
$ws.Range("N7").Value = 0.1806851486826251
$ws.Range("O7").Value = 0.6953460040961424
$ws.Range("P7").Value = 20.18527741693543

# Sheet 12
$ws = $wb.Worksheets.Item(12)
$ws.Range("D7").Value = 0.4433215424468767
$ws.Range("E7").Value = 0.5708552006732649
$ws.Range("F7").Value = 0.2121279089469113
$ws.Range("G7").Value = 0.6658239575494987
$ws.Range("H7").Value = 25.00652885694351
$ws.Range("L7").Value = 0.4168784746538839
$ws.Range("M7").Value = 0.5472476230423133
$ws.Range("N7").Value = 0.2021769619747467
$ws.Range("O7").Value = 0.6456612692843546
$ws.Range("P7").Value = 23.67559552856849

# Sheet 13
$ws = $wb.Worksheets.Item(13)
$ws.Range("D7").Value = 0.01926107829590986
$ws.Range("E7").Value = 0.09942345903022236
$ws.Range("F7").Value = 0.1338056038767039
$ws.Range("G7").Value = 0.1387842869200612
$ws.Range("H7").Value = 14.71414692364147
$ws.Range("L7").Value = 0.01958623205003894
$ws.Range("M7").Value = 0.1038061133648449
$ws.Range("N7").Value = 0.1396715610001706
$ws.Range("O7").Value = 0.1399508201120627
$ws.Range("P7").Value = 15.18774924299386

# Sheet 14
$ws = $wb.Worksheets.Item(14)
$ws.Range("D7").Value = 0.01628485489803137
$ws.Range("E7").Value = 0.09496478459446893
$ws.Range("F7").Value = 0.1540408183477182
$ws.Range("G7").Value = 0.1276121267671352
$ws.Range("H7").Value = 15.55192190762534
$ws.Range("L7").Value = 0.01632941127840756
$ws.Range("M7").Value = 0.09475106205950076
$ws.Range("N7").Value = 0.1567278526708343
$ws.Range("O7").Value = 0.1277865848921848
$ws.Range("P7").Value = 15.59198430293306

# Sheet 15
$ws = $wb.Worksheets.Item(15)
$ws.Range("D7").Value = 0.03039417414469136
$ws.Range("E7").Value = 0.1249104029861712
$ws.Range("F7").Value = 0.2523922033975708
$ws.Range("G7").Value = 0.17433925015524
$ws.Range("H7").Value = 22.04587138245413
$ws.Range("L7").Value = 0.03036237131563837
$ws.Range("M7").Value = 0.1247897074141153
$ws.Range("N7").Value = 0.252117800290957
$ws.Range("O7").Value = 0.1742480166763409
$ws.Range("P7").Value = 22.03087689554467

# Sheet 16
$ws = $wb.Worksheets.Item(16)
$ws.Range("D7").Value = 0.2113255131626701
$ws.Range("E7").Value = 0.373130555570257
$ws.Range("F7").Value = 0.1329607055997786
$ws.Range("G7").Value = 0.4597015479228563
$ws.Range("H7").Value = 14.1657266600943
$ws.Range("L7").Value = 0.1928721468812845
$ws.Range("M7").Value = 0.3495385903618834
$ws.Range("N7").Value = 0.1244362025035707
$ws.Range("O7").Value = 0.439172115327561
$ws.Range("P7").Value = 13.53033922751493

# Sheet 17
$ws = $wb.Worksheets.Item(17)
$ws.Range("D7").Value = 0.223171839511746
$ws.Range("E7").Value = 0.371367614300916
$ws.Range("F7").Value = 0.1413642043835387
$ws.Range("G7").Value = 0.472410668287398
$ws.Range("H7").Value = 15.08583435060929
$ws.Range("L7").Value = 0.20199957143429
$ws.Range("M7").Value = 0.3324475471469159
$ws.Range("N7").Value = 0.1263587920573767
$ws.Range("O7").Value = 0.4494436243115369
$ws.Range("P7").Value = 13.50673027190282

# Sheet 18
$ws = $wb.Worksheets.Item(18)
$ws.Range("D7").Value = 0.2229140131546724
$ws.Range("E7").Value = 0.3541756080780334
$ws.Range("F7").Value = 0.151838766880737
$ws.Range("G7").Value = 0.4721377057116625
$ws.Range("H7").Value = 15.1892837235934
$ws.Range("L7").Value = 0.3329011307212459
$ws.Range("M7").Value = 0.3918314459511077
$ws.Range("N7").Value = 0.1608660613368267
$ws.Range("O7").Value = 0.5769758493396807
$ws.Range("P7").Value = 17.08925434764643

# Sheet 19
$ws = $wb.Worksheets.Item(19)
$ws.Range("D7").Value = 0.1239779984273855
$ws.Range("E7").Value = 0.2700361536234256
$ws.Range("F7").Value = 0.3857741693797568
$ws.Range("G7").Value = 0.352105095713461
$ws.Range("H7").Value = 48.84502021247369

# Sheet 20
$ws = $wb.Worksheets.Item(20)
$ws.Range("D7").Value = 0.2259296198691152
$ws.Range("E7").Value = 0.3994213745820632
$ws.Range("F7").Value = 0.6949309772357419
$ws.Range("G7").Value = 0.4753205443373085
$ws.Range("H7").Value = 61.01509009595573

# Sheet 21
$ws = $wb.Worksheets.Item(21)
$ws.Range("D7").Value = 0.1181884301055577
$ws.Range("E7").Value = 0.263682557078448
$ws.Range("F7").Value = 0.6261524469939328
$ws.Range("G7").Value = 0.3437854419628
$ws.Range("H7").Value = 45.79182914756211

# Sheet 22
$ws = $wb.Worksheets.Item(22)
$ws.Range("D7").Value = 3.550906062299529
$ws.Range("E7").Value = 1.304004217518093
$ws.Range("F7").Value = 0.4800791138254637
$ws.Range("G7").Value = 1.884384796770428
$ws.Range("H7").Value = 63.70369238720948

# Sheet 23
$ws = $wb.Worksheets.Item(23)
$ws.Range("D7").Value = 6.929291472209757
$ws.Range("E7").Value = 2.145603192834613
$ws.Range("F7").Value = 0.8229969750038781
$ws.Range("G7").Value = 2.632354739052044
$ws.Range("H7").Value = 91.5653413564839

# Sheet 24
$ws = $wb.Worksheets.Item(24)
$ws.Range("D7").Value = 10.8467900123392
$ws.Range("E7").Value = 2.660979836234757
$ws.Range("F7").Value = 1.23685460830496
$ws.Range("G7").Value = 3.29344652489443
$ws.Range("H7").Value = 91.3061774051198
